# Add 2 new scenes to the Scene table (rows 21 and 22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 21: 落潮小径 (riverside trail scene)
$ws.Range("A21").Value = 13010017
$ws.Range("B21").Value = "落潮小径"
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 1038
$ws.Range("E21").Value = 459
$ws.Range("H21").Value = 18
$ws.Range("I21").Value = "default"

# New row 22: 月光林地 (moonlit grove scene)
$ws.Range("A22").Value = 13010018
$ws.Range("B22").Value = "月光林地"
$ws.Range("C22").Value = 31
$ws.Range("D22").Value = 723
$ws.Range("E22").Value = 327
$ws.Range("H22").Value = 19
$ws.Range("I22").Value = "default"

# Expand the table (ListObject) to include the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I22"))

# Match the final selection from the authored change
$ws.Range("H22").Select()
